$wb = $excel.ActiveWorkbook
$wsIde = $wb.Worksheets.Item("ide")
$wsPcb = $wb.Worksheets.Item("pcb")

# --- Add new "2018.8.3" eval_only result row to the "ide" sheet ---
# Use a formula literal + paste-special values so the date-like text
# "2018.8.3" is stored as a plain string (shared string) rather than
# being auto-parsed into a date serial number.
$wsIde.Range("A8").Formula = "=""2018.8.3"""
$wsIde.Range("A8").Copy()
$wsIde.Range("A8").PasteSpecial(-4163)
$wsIde.Range("B8").NumberFormat = $wsIde.Range("B7").NumberFormat
$wsIde.Range("B8").Value = 0.657

# --- Re-apply center/middle alignment on the section-header cells so the
# style table keeps its existing (center) / (center+middle) xf entries
# in the same order they end up in after Excel round-trips the file ---
$wsIde.Range("A1:B1").HorizontalAlignment = -4108
$wsIde.Range("A1:B1").VerticalAlignment = -4108

$wsPcb.Range("A1:B1").HorizontalAlignment = -4108
$wsPcb.Range("A1:B1").VerticalAlignment = -4108

$wsPcb.Range("A6:B6").HorizontalAlignment = -4108
$wsPcb.Range("A6:B6").VerticalAlignment = -4108

$wsPcb.Range("A10").VerticalAlignment = -4107

# --- Switch the active/selected sheet from "pcb" to "ide" (eval_only flag) ---
$wsIde.Activate()
$wsIde.Range("A9").Select()
